$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.560.06'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.19%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.506.10'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.48%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '606.26'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.02%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '151.84'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.34%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.504.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.51%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.487'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.12%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.143'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.61%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.64'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +7.38%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.433'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.82%  '

$ws.Range('E13').Value = '  -1.76%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.27'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.64%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.101.11'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.47%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.508.05'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.45%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.487.15'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.10%  '

$ws.Range('E18').Value = '  -0.66%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.51'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.85%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.47'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.18%  '

$ws.Range('E21').Value = '  +2.69%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '446.82'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.29%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.629'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.79%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '78.32'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.08%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.647.48'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.46%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000127'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.58%  '

$ws.Range('E27').Value = '  -0.01%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.74'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.52%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.04'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.60%  '

$ws.Range('E30').Value = '  +0.16%  '

$ws.Range('E31').Value = '  +5.36%  '

$ws.Range('E32').Value = '  +5.57%  '

$ws.Range('E33').Value = '  -0.02%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.61'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.90%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.14'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.33%  '

$ws.Range('E36').Value = '  +0.89%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.497.48'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.47%  '

$ws.Range('E38').Value = '  -0.13%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.30'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.88%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '179.69'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.55%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.04%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0896'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.50%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.44'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.38%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.892'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.47%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '30.23'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +7.26%  '

$ws.Range('E47').Value = '  +2.96%  '

$ws.Range('E48').Value = '  +4.21%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.54'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.66%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.61'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.24%  '

$ws.Range('E51').Value = '  +1.84%  '

